# Apply the "stock.yaml completed" update to the ABCAPITAL.NS 1wk sheet:
#  - add a new "backup" column (R) with a header + mostly-zero values
#  - correct a handful of stale "detect_structure" (Q) / "isPivot" (O) values
#  - append 14 new weekly rows (356-369) of OHLCV + derived data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column R: header "backup" (copy the bold/bordered header style
#    from A1), then the data rows default to 0.
# ---------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("R1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("R1").Value = "backup"

$ws.Range("R2:R355").Value = 0

# A few rows keep a non-zero "backup" value.
$backupOverrides = @{
    132 = 2
    156 = 2
    244 = 1
    259 = 1
    271 = 2
    283 = 2
}
foreach ($r in $backupOverrides.Keys) {
    $ws.Cells.Item($r, 18).Value = $backupOverrides[$r]
}

# ---------------------------------------------------------------------
# 2. A handful of previously-mis-detected rows get their Q
#    ("detect_structure") value reset to 0.
# ---------------------------------------------------------------------
$ws.Range("Q44").Value = 0
$ws.Range("Q54").Value = 0
$ws.Range("Q57").Value = 0

# Row 355's "isPivot" (O) flag is corrected from 0 to 2.
$ws.Range("O355").Value = 2

# ---------------------------------------------------------------------
# 3. Append the 14 new weekly rows (356-369). Column A keeps the same
#    date/time number format as the rest of the sheet.
# ---------------------------------------------------------------------
$ws.Range("A355").Copy() | Out-Null
$ws.Range("A356:A369").PasteSpecial(-4122) | Out-Null

$newRows = @(
    @(356,45453,233.3699951171875,243.3000030517578,229,239.9600067138672,239.9600067138672,37955390,2024,6,10,0,0,0,24,0,0,2),
    @(357,45460,239.9600067138672,246.8999938964844,225,241.3000030517578,241.3000030517578,39398324,2024,6,17,0,0,0,25,1,0,0),
    @(358,45467,239.7100067138672,243.3899993896484,232.7599945068359,238.7100067138672,238.7100067138672,23368740,2024,6,24,0,0,0,26,0,0,0),
    @(359,45474,239.3999938964844,242.4100036621094,234.3999938964844,236.0099945068359,236.0099945068359,16168145,2024,7,1,0,0,0,27,0,0,0),
    @(360,45481,236.0099945068359,237.0399932861328,224.2899932861328,224.75,224.75,15669754,2024,7,8,0,0,0,28,0,0,0),
    @(361,45488,225.3399963378906,228.1799926757812,215.0500030517578,215.6399993896484,215.6399993896484,15606722,2024,7,15,0,0,0,29,0,0,1),
    @(362,45495,215,224.4499969482422,205.5299987792969,223.4299926757812,223.4299926757812,22179365,2024,7,22,0,0,0,30,0,0,0),
    @(363,45502,225,229.3999938964844,208,211.5899963378906,211.5899963378906,34437693,2024,7,29,0,0,0,31,0,0,2),
    @(364,45509,204,214.1499938964844,200.0500030517578,211.8000030517578,211.8000030517578,28833609,2024,8,5,0,0,0,32,2,0,0),
    @(365,45516,211.6000061035156,216.9799957275391,201.6000061035156,216.0899963378906,216.0899963378906,24676351,2024,8,12,0,0,0,33,0,0,0),
    @(366,45523,217.1199951171875,226,211.9199981689453,219.0800018310547,219.0800018310547,29512913,2024,8,19,0,0,0,34,0,0,0),
    @(367,45530,220.5,225.6999969482422,218.0899963378906,220.1100006103516,220.1100006103516,25703323,2024,8,26,0,0,0,35,0,0,0),
    @(368,45537,221,229.4499969482422,215.6399993896484,216.7799987792969,216.7799987792969,28590380,2024,9,2,0,0,0,36,0,0,0),
    @(369,45544,217.8999938964844,227.1799926757812,212.8000030517578,226.1000061035156,226.1000061035156,18231324,2024,9,9,0,0,0,37,0,0,0)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $vals = New-Object 'object[,]' 1,17
    for ($i = 0; $i -lt 17; $i++) {
        $vals[0,$i] = $row[$i + 1]
    }
    $ws.Range("A$r`:Q$r").Value = $vals
}

Write-Output "edit applied"
